# Registro de aportes y adecuar miembros a 2024
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Ingreso": correct the last existing row amount and append new
# contribution rows (aportes) for late Dec 2023 / early Jan 2024.
# ---------------------------------------------------------------------
$wsIngreso = $wb.Worksheets.Item("Ingreso")

# Row 576 amount correction: 750 -> 350
$wsIngreso.Cells.Item(576, 3).Value = 350

# New rows 577-581
$wsIngreso.Cells.Item(577, 1).Value = 45284
$wsIngreso.Cells.Item(577, 2).Value = "Carlos"
$wsIngreso.Cells.Item(577, 3).Value = 500
$wsIngreso.Cells.Item(577, 4).Value = "Aporte"

$wsIngreso.Cells.Item(578, 1).Value = 45284
$wsIngreso.Cells.Item(578, 2).Value = "Omaury"
$wsIngreso.Cells.Item(578, 3).Value = 100
$wsIngreso.Cells.Item(578, 4).Value = "Aporte"

$wsIngreso.Cells.Item(579, 1).Value = 45284
$wsIngreso.Cells.Item(579, 2).Value = "Alfredo"
$wsIngreso.Cells.Item(579, 3).Value = 100
$wsIngreso.Cells.Item(579, 4).Value = "Aporte"

$wsIngreso.Cells.Item(580, 1).Value = 45284
$wsIngreso.Cells.Item(580, 2).Value = "Julio"
$wsIngreso.Cells.Item(580, 3).Value = 100
$wsIngreso.Cells.Item(580, 4).Value = "Aporte"

$wsIngreso.Cells.Item(581, 1).Value = 45292
$wsIngreso.Cells.Item(581, 2).Value = "Johan"
$wsIngreso.Cells.Item(581, 3).Value = 300
$wsIngreso.Cells.Item(581, 4).Value = "Aporte"

# ---------------------------------------------------------------------
# Sheet "Gastos": append a new expense row for the referee/water/ice.
# ---------------------------------------------------------------------
$wsGastos = $wb.Worksheets.Item("Gastos")

$wsGastos.Cells.Item(75, 1).Value = 45284
$wsGastos.Cells.Item(75, 2).Value = "Arbitro, agua y hielo"
$wsGastos.Cells.Item(75, 3).Value = 1200

# ---------------------------------------------------------------------
# Sheet "Cuentas por cobrar": log a new "Cobros" entry with a formula.
# Column A has no sheet-wide style, so copy the date format from the
# cell above (A5) instead of assigning a NumberFormat string, which
# would otherwise register as a brand new (duplicate) style.
# ---------------------------------------------------------------------
$wsCobrar = $wb.Worksheets.Item("Cuentas por cobrar")

$wsCobrar.Cells.Item(5, 1).Copy()
$wsCobrar.Cells.Item(6, 1).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$wsCobrar.Cells.Item(6, 1).Value = 45270
$wsCobrar.Cells.Item(6, 2).Value = "Joel"
$wsCobrar.Cells.Item(6, 3).Value = "Cobros"
$wsCobrar.Cells.Item(6, 4).Formula = "=500+50+100+60-100"

# ---------------------------------------------------------------------
# Selection / view bookkeeping to mirror the author's final state.
# "Ingreso" is the tab left active/selected, so it must be the last
# sheet touched here.
# ---------------------------------------------------------------------
$wsActividad = $wb.Worksheets.Item("Actividad GOAT")

[void]$wsGastos.Range("A75").Select()
[void]$wsCobrar.Range("A7").Select()
[void]$wsActividad.Range("D35").Select()
$wsIngreso.Activate()
[void]$wsIngreso.Range("C581").Select()
